$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.943.54"
$ws.Range("E2").Value = "'  -0.35%  "
$ws.Range("D3").Value = "'2.418.28"
$ws.Range("E3").Value = "'  -0.09%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'561.88"
$ws.Range("E5").Value = "'  -0.24%  "
$ws.Range("D6").Value = "'142.76"
$ws.Range("E6").Value = "'  -0.93%  "
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E8").Value = "'  -0.61%  "
$ws.Range("E9").Value = "'  -0.19%  "
$ws.Range("E10").Value = "'  -0.80%  "
$ws.Range("E11").Value = "'  -4.15%  "
$ws.Range("E12").Value = "'  -1.81%  "
$ws.Range("D13").Value = "'26.13"
$ws.Range("E13").Value = "'  +0.32%  "
$ws.Range("E14").Value = "'  -2.25%  "
$ws.Range("D15").Value = "'2.838.73"
$ws.Range("E15").Value = "'  -0.65%  "
$ws.Range("D16").Value = "'61.884.36"
$ws.Range("E16").Value = "'  -0.27%  "
$ws.Range("D17").Value = "'2.419.29"
$ws.Range("E17").Value = "'  -0.03%  "
$ws.Range("D18").Value = "'11.30"
$ws.Range("E18").Value = "'  +0.39%  "
$ws.Range("D19").Value = "'323.12"
$ws.Range("E19").Value = "'  -0.44%  "
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.82"
$ws.Range("E20").Value = "'  +0.96%  "
$ws.Range("B21").Value = "'Polkadot"
$ws.Range("C21").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'4.13"
$ws.Range("E21").Value = "'  -1.61%  "
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("D23").Value = "'66.92"
$ws.Range("E23").Value = "'  +2.41%  "
$ws.Range("D24").Value = "'1.74"
$ws.Range("E24").Value = "'  +0.11%  "
$ws.Range("D25").Value = "'8.75"
$ws.Range("E25").Value = "'  -3.55%  "
$ws.Range("D26").Value = "'553.19"
$ws.Range("E26").Value = "'  -6.04%  "
$ws.Range("D27").Value = "'2.537.63"
$ws.Range("E27").Value = "'  +0.42%  "
$ws.Range("E28").Value = "'  +0.16%  "
$ws.Range("D29").Value = "'0.0₃0929"
$ws.Range("E29").Value = "'  -1.58%  "
$ws.Range("D30").Value = "'8.18"
$ws.Range("E30").Value = "'  -0.85%  "
$ws.Range("D31").Value = "'1.38"
$ws.Range("E31").Value = "'  -4.83%  "
$ws.Range("E32").Value = "'  -1.96%  "
$ws.Range("E33").Value = "'  -0.95%  "
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "'  -3.98%  "
$ws.Range("E35").Value = "'  -0.04%  "
$ws.Range("D36").Value = "'4.72"
$ws.Range("E36").Value = "'  -1.15%  "
$ws.Range("E37").Value = "'  -1.71%  "
$ws.Range("B38").Value = "'Monero"
$ws.Range("C38").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'152.38"
$ws.Range("E38").Value = "'  -1.08%  "
$ws.Range("B39").Value = "'RenderToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'5.43"
$ws.Range("E39").Value = "'  -5.21%  "
$ws.Range("D40").Value = "'18.61"
$ws.Range("E40").Value = "'  -0.32%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "'  -1.14%  "
$ws.Range("D42").Value = "'0.992"
$ws.Range("E42").Value = "'  -0.71%  "
$ws.Range("D43").Value = "'147.00"
$ws.Range("E43").Value = "'  -2.35%  "
$ws.Range("D44").Value = "'2.22"
$ws.Range("E44").Value = "'  -4.98%  "
$ws.Range("D45").Value = "'3.63"
$ws.Range("E45").Value = "'  -0.59%  "
$ws.Range("D46").Value = "'0.0526"
$ws.Range("E46").Value = "'  -2.29%  "
$ws.Range("E47").Value = "'  +0.38%  "
$ws.Range("D48").Value = "'19.79"
$ws.Range("E48").Value = "'  -2.73%  "
$ws.Range("E49").Value = "'  -0.65%  "
$ws.Range("D50").Value = "'0.0227"
$ws.Range("E50").Value = "'  -0.72%  "
